$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '44.220.52'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +1.21%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.255.00'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +0.41%  '

$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '307.99'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -4.56%  '

$ws.Range('E6').Value = '  -2.26%  '

$ws.Range('E7').Value = '  -0.63%  '

$ws.Range('E8').Value = '  +0.09%  '

$ws.Range('E9').Value = '  -3.32%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.59'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -4.54%  '

$ws.Range('E11').Value = '  -0.92%  '

$ws.Range('E12').Value = '  -4.70%  '

$ws.Range('E13').Value = '  -1.71%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.597.07'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.33%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.289.72'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +1.99%  '

$ws.Range('E16').Value = '  -1.98%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.85'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -2.63%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '44.068.63'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +1.03%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.94'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -5.36%  '

$ws.Range('E20').Value = '  -0.74%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.35'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -3.20%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '65.56'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.53%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '243.29'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +2.87%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.96'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -6.48%  '

$ws.Range('E25').Value = '  -8.81%  '

$ws.Range('E26').Value = '  +0.01%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.15'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.03%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.15'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -1.92%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '36.84'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.25%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.20'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -1.76%  '

$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '20.16'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.06%  '

$ws.Range('B32').Value = 'LidoDAOToken'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.59'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +14.07%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '157.48'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.86%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0829'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -2.94%  '

$ws.Range('E35').Value = '  -1.11%  '

$ws.Range('E36').Value = '  -0.52%  '

$ws.Range('E37').Value = '  -4.38%  '

$ws.Range('E38').Value = '  -3.88%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '15.33'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.74%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.91'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -8.03%  '

$ws.Range('E41').Value = '  -9.92%  '

$ws.Range('E42').Value = '  -3.40%  '

$ws.Range('E43').Value = '  +0.08%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.765.68'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -2.18%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '88.75'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +7.73%  '

$ws.Range('E46').Value = '  -0.51%  '

$ws.Range('E47').Value = '  -3.28%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '101.76'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -1.41%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.28'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -2.03%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '70.52'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -5.37%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '55.53'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -5.35%  '
